$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated price (D) and volume (E) values to match the latest crypto data snapshot.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '22.410.89'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.563.17'
$ws.Range('D3').Style = "Normal"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '285.40'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3634'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.86%  '
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3336'
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.128'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07409'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.81'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.932'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.893'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.564.30'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('E18').Value = '  -3.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06689'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.355'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.12'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '22.397.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.411'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.41%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.558'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '149.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.41'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.990'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.739.57'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.060'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.139'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.998'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.814'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08240'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.05%  '
$ws.Range('E37').Value = '  -2.62%  '
$ws.Range('E38').Value = '  -5.87%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06386'
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2209'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.333'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.17'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6085'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5759'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('E48').Value = '  -3.28%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.53'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.214'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07210'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.50%  '
